$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the header text in G4 from "Estado" to "¿Terminada?"
$ws.Range("G4").Value = "¿Terminada?"

# Update the selected cell to G4 (matches the author's last selection in the diff)
$ws.Range("G4").Select()
